# Small adjustments to solar PV DK.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Unhide DKW1/DKE1 2040 "Distributed Energy" rows (38-41) and update values ---
$ws.Rows.Item(38).Hidden = $false
$ws.Range("F38").Value = 10200

$ws.Rows.Item(39).Hidden = $false
$ws.Range("F39").Value = 2700

$ws.Rows.Item(40).Hidden = $false
$ws.Range("F40").Value = 7300

$ws.Rows.Item(41).Hidden = $false
$ws.Range("F41").Value = 2700

# --- Hide rows 62-71 (2030 Solar PV / Onshore Wind rows) ---
for ($r = 62; $r -le 71; $r++) {
    $ws.Rows.Item($r).Hidden = $true
}

# --- Unhide rows 81-86 (DKE1/DKW1 2030 rows) ---
for ($r = 81; $r -le 86; $r++) {
    $ws.Rows.Item($r).Hidden = $false
}

# --- Hide row 90 ---
$ws.Rows.Item(90).Hidden = $true

# --- Add new rows 91-92 for DKE1/DKW1 2040 Solar PV Distributed Energy ---
$ws.Range("A91").Value = "DKE1"
$ws.Range("C91").Value = "Solar PV"
$ws.Range("D91").Value = "Distributed Energy"
$ws.Range("E91").Value = 2040
$ws.Range("F91").Value = 2400

$ws.Range("A92").Value = "DKW1"
$ws.Range("C92").Value = "Solar PV"
$ws.Range("D92").Value = "Distributed Energy"
$ws.Range("E92").Value = 2040
$ws.Range("F92").Value = 5600

# --- Update autofilter: Node (col 0) = DKE1/DKW1; Year (col 4) = blank/2030/2040 ---
$nodeVals = @("DKE1", "DKW1")
$ws.Range("A1:J90").AutoFilter(1, $nodeVals, 7)

$yearVals = @("", "2030", "2040")
$ws.Range("A1:J90").AutoFilter(5, $yearVals, 7)

# --- Update selection to match the recorded cursor position ---
$ws.Range("D96").Select()
